# A new entry was logged on the "2024" sheet (most-recent-first list of
# September notification entries). Insert a new row above the current
# top data row (row 29), which shifts every row below it (29 -> 30, ...,
# 55 -> 56) down by one and grows the sheet dimension to A1:Y56.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2024")

$ws.Rows.Item(29).Insert()

# Populate the newly inserted row with the new "September_Details" /
# "September_Date" entry (columns R and S); the remaining columns on
# this row stay blank, same as every other log row.
$ws.Range("R29").Value = "exclusive on axis"
$ws.Range("S29").Value = "2024-09-04 13:21:05"
